$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; this shifts old rows 25-56 down to 26-57
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with this week's data (matches format of surrounding rows)
$ws.Range("A25").Value = 8
$ws.Range("B25").Value = "Terminal La Palmera de La Serena"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 44771
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 100114007
$ws.Range("G25").Value = "Jengibre"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 480
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("N25").Value = "$/caja 13 kilos"
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 1115
$ws.Range("Q25").Value = 13
$ws.Range("R25").Value = "Hortaliza"
